$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "r"
$ws.Range("D7").Value = "hth"
$ws.Range("E9").Value = "ht"
$ws.Range("E9").Select()
